$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: read the "ExpectedRate" (column D) cell as the currency-formatted
# text that "ActualRate" (column E) should now mirror. Some D cells already
# hold literal text (e.g. "$77.91"), others hold a real number formatted as
# currency - normalise both to the same "$#,##0.00" text representation.
function Get-CurrencyText($cell) {
    $v = $cell.Value2
    if ($v -is [string]) {
        return $v
    } else {
        return [string]::Format("{0:C2}", $v)
    }
}

# Rows whose FedEx rate verification now passes: ActualRate is corrected to
# equal ExpectedRate, and Result flips from FAIL to PASS.
$rows = 2,3,4,5,6,7,8,9,10,11,12,13,14,15,30

foreach ($r in $rows) {
    $dCell = $ws.Range("D$r")
    $eCell = $ws.Range("E$r")
    $txt = Get-CurrencyText $dCell

    # Force the assignment to stay plain text (matching the existing
    # ActualRate cells) instead of letting Excel auto-convert a
    # "$"-looking string into a formatted number.
    $eCell.NumberFormat = "@"
    $eCell.Value = $txt
    $eCell.Style = "Normal"

    $ws.Range("F$r").Value = "PASS"
}

# Touch a cell in column H so the sheet's used range / dimension extends
# to include it, matching the recorded change (A1:G40 -> A1:H40).
$ws.Range("H1").Font.Bold = $ws.Range("H1").Font.Bold
